# Allocation rule summary tables updated with 5 mi and 10 mi radius columns
# (adds column F "Within 5 mile(s)..." and column G "Within 10 mile(s)..."
#  to both the Means and Standard Deviations sheets, and refreshes the
#  "Total Cancer Risk" / "Total Respiratory" rows for the existing columns).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: Means ---
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"
$ws1.Range("F2").Value = 51
$ws1.Range("G2").Value = 70
$ws1.Range("F3").Value = 43
$ws1.Range("G3").Value = 24
$ws1.Range("F4").Value = 6.1
$ws1.Range("G4").Value = 5.7
$ws1.Range("F5").Value = 4.5
$ws1.Range("G5").Value = 5.5
$ws1.Range("F6").Value = 37
$ws1.Range("G6").Value = 51
$ws1.Range("F7").Value = 14
$ws1.Range("G7").Value = 9.7
$ws1.Range("F8").Value = 12
$ws1.Range("G8").Value = 8
$ws1.Range("B9").Value = 29
$ws1.Range("C9").Value = 29
$ws1.Range("D9").Value = 30
$ws1.Range("E9").Value = 30
$ws1.Range("F9").Value = 30
$ws1.Range("G9").Value = 30
$ws1.Range("B10").Value = 0.37
$ws1.Range("C10").Value = 0.36
$ws1.Range("D10").Value = 0.37
$ws1.Range("E10").Value = 0.4
$ws1.Range("F10").Value = 0.43
$ws1.Range("G10").Value = 0.41

# --- Sheet 2: Standard Deviations ---
$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"
$ws2.Range("F2").Value = 32
$ws2.Range("G2").Value = 30
$ws2.Range("F3").Value = 33
$ws2.Range("G3").Value = 30
$ws2.Range("F4").Value = 6.7
$ws2.Range("G4").Value = 6.3
$ws2.Range("F5").Value = 6.8
$ws2.Range("G5").Value = 8.4
$ws2.Range("F6").Value = 16
$ws2.Range("G6").Value = 25
$ws2.Range("F7").Value = 11
$ws2.Range("G7").Value = 10
$ws2.Range("F8").Value = 12
$ws2.Range("G8").Value = 10
$ws2.Range("B9").Value = 10
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 1.8
$ws2.Range("G9").Value = 1.1
$ws2.Range("B10").Value = 0.14
$ws2.Range("C10").Value = 0.086
$ws2.Range("D10").Value = 0.052
$ws2.Range("E10").Value = 0.032
$ws2.Range("F10").Value = 0.05
$ws2.Range("G10").Value = 0.048
